$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 272.25
$ws.Range("I15").Value = 272.25
$ws.Range("K15").Value = 816.75
$ws.Range("M15").Value = -647.75
$ws.Range("H62").Value = 1916.3226
$ws.Range("I62").Value = 2022.1666
$ws.Range("J62").Value = 1769.7693
$ws.Range("K62").Value = 2022.1666
$ws.Range("L62").Value = 1769.7693
$ws.Range("M62").Value = -1398.1666
$ws.Range("N62").Value = -3017.7693
$ws.Range("H65").Value = 1916.3226
$ws.Range("I65").Value = 2022.1666
$ws.Range("J65").Value = 1769.7693
$ws.Range("K65").Value = 10110.833
$ws.Range("L65").Value = 8848.8465
$ws.Range("M65").Value = -6990.833000000001
$ws.Range("N65").Value = -15088.8465
$ws.Range("H80").Value = 5962.3887
$ws.Range("I80").Value = 452.45456
$ws.Range("K80").Value = 1357.36368
$ws.Range("M80").Value = -359.3636799999999
$ws.Range("H83").Value = 5962.3887
$ws.Range("I83").Value = 452.45456
$ws.Range("K83").Value = 4072.09104
$ws.Range("M83").Value = 919.9089599999998
$ws.Range("H132").Value = 1951.4
$ws.Range("I132").Value = 1781.1724
$ws.Range("K132").Value = 5343.5172
$ws.Range("M132").Value = -2813.5172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 730.6842
$ws.Range("I2").Value = 591.19354
$ws.Range("K2").Value = 591.19354
$ws.Range("M2").Value = -478.19354
$ws.Range("H116").Value = 730.6842
$ws.Range("I116").Value = 591.19354
$ws.Range("K116").Value = 591.19354
$ws.Range("M116").Value = 1702.80646
$ws.Range("H132").Value = 1268061.9
$ws.Range("I132").Value = 1432.6666
$ws.Range("J132").Value = 6255414.5
$ws.Range("K132").Value = 4297.9998
$ws.Range("L132").Value = 18766243.5
$ws.Range("M132").Value = -1767.9998
$ws.Range("N132").Value = -18771303.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 730.6842
$ws.Range("I3").Value = 591.19354
$ws.Range("K3").Value = 591.19354
$ws.Range("M3").Value = -477.19354
$ws.Range("H20").Value = 19556.285
$ws.Range("I20").Value = 1853.4546
$ws.Range("J20").Value = 84466.664
$ws.Range("K20").Value = 1853.4546
$ws.Range("L20").Value = 84466.664
$ws.Range("M20").Value = -1606.4546
$ws.Range("N20").Value = -84960.664
$ws.Range("H29").Value = 4471
$ws.Range("I29").Value = 3783.75
$ws.Range("K29").Value = 3783.75
$ws.Range("M29").Value = -3494.75
$ws.Range("H134").Value = 33563.145
$ws.Range("I134").Value = 5987.84
$ws.Range("J134").Value = 102501.4
$ws.Range("K134").Value = 17963.52
$ws.Range("L134").Value = 307504.2
$ws.Range("M134").Value = -15428.52
$ws.Range("N134").Value = -312574.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4043.1868
$ws.Range("I31").Value = 1814.5227
$ws.Range("J31").Value = 7206.4517
$ws.Range("K31").Value = 1814.5227
$ws.Range("L31").Value = 7206.4517
$ws.Range("M31").Value = -1519.5227
$ws.Range("N31").Value = -7796.4517
$ws.Range("H34").Value = 4043.1868
$ws.Range("I34").Value = 1814.5227
$ws.Range("J34").Value = 7206.4517
$ws.Range("K34").Value = 1814.5227
$ws.Range("L34").Value = 7206.4517
$ws.Range("M34").Value = -1612.5227
$ws.Range("N34").Value = -7610.4517
$ws.Range("H58").Value = 334793.3
$ws.Range("I58").Value = 1299.3846
$ws.Range("J58").Value = 2502503.8
$ws.Range("K58").Value = 1299.3846
$ws.Range("L58").Value = 2502503.8
$ws.Range("M58").Value = -1096.3846
$ws.Range("N58").Value = -2502909.8
$ws.Range("H136").Value = 334793.3
$ws.Range("I136").Value = 1299.3846
$ws.Range("J136").Value = 2502503.8
$ws.Range("K136").Value = 3898.1538
$ws.Range("L136").Value = 7507511.399999999
$ws.Range("M136").Value = -1348.1538
$ws.Range("N136").Value = -7512611.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1282.921
$ws.Range("I5").Value = 449.7143
$ws.Range("J5").Value = 1768.9584
$ws.Range("K5").Value = 1349.1429
$ws.Range("L5").Value = 5306.8752
$ws.Range("M5").Value = -1237.1429
$ws.Range("N5").Value = -5530.8752
$ws.Range("H12").Value = 3258210.2
$ws.Range("J12").Value = 50138.7
$ws.Range("L12").Value = 150416.1
$ws.Range("N12").Value = -150762.1
$ws.Range("H20").Value = 5352.4287
$ws.Range("J20").Value = 6161.1665
$ws.Range("L20").Value = 18483.4995
$ws.Range("N20").Value = -18937.4995
$ws.Range("H22").Value = 2400.1667
$ws.Range("J22").Value = 2800
$ws.Range("L22").Value = 8400
$ws.Range("N22").Value = -8738
$ws.Range("H27").Value = 2400.1667
$ws.Range("J27").Value = 2800
$ws.Range("L27").Value = 8400
$ws.Range("N27").Value = -8604
$ws.Range("H135").Value = 1282.921
$ws.Range("I135").Value = 449.7143
$ws.Range("J135").Value = 1768.9584
$ws.Range("K135").Value = 4047.4287
$ws.Range("L135").Value = 15920.6256
$ws.Range("M135").Value = -1512.4287
$ws.Range("N135").Value = -20990.6256

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10668.667
$ws.Range("I70").Value = 10668.667
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 10668.667
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -10398.667
$ws.Range("H73").Value = 10668.667
$ws.Range("I73").Value = 10668.667
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 10668.667
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -9732.666999999999
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 40650.73
$ws.Range("I7").Value = 57539.11
$ws.Range("J7").Value = 2651.875
$ws.Range("K7").Value = 57539.11
$ws.Range("L7").Value = 2651.875
$ws.Range("M7").Value = -57427.11
$ws.Range("N7").Value = -2875.875
$ws.Range("H16").Value = 1020.94116
$ws.Range("I16").Value = 1072.75
$ws.Range("J16").Value = 192
$ws.Range("K16").Value = 1072.75
$ws.Range("L16").Value = 192
$ws.Range("M16").Value = -902.75
$ws.Range("N16").Value = -532
$ws.Range("H122").Value = 1897674.6
$ws.Range("I122").Value = 2555826
$ws.Range("J122").Value = 669125.25
$ws.Range("K122").Value = 7667478
$ws.Range("L122").Value = 2007375.75
$ws.Range("M122").Value = -7665028
$ws.Range("N122").Value = -2012275.75
$ws.Range("H126").Value = 40650.73
$ws.Range("I126").Value = 57539.11
$ws.Range("J126").Value = 2651.875
$ws.Range("K126").Value = 172617.33
$ws.Range("L126").Value = 7955.625
$ws.Range("M126").Value = -170147.33
$ws.Range("N126").Value = -12895.625
$ws.Range("H133").Value = 110163
$ws.Range("J133").Value = 110163
$ws.Range("L133").Value = 110163
$ws.Range("N133").Value = -115223

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2620.4
$ws.Range("I4").Value = 2002
$ws.Range("J4").Value = 2775
$ws.Range("K4").Value = 2002
$ws.Range("L4").Value = 2775
$ws.Range("M4").Value = -1889
$ws.Range("N4").Value = -3001
$ws.Range("H126").Value = 528.4838999999999
$ws.Range("I126").Value = 446
$ws.Range("J126").Value = 1298.3334
$ws.Range("K126").Value = 1338
$ws.Range("L126").Value = 3895.0002
$ws.Range("M126").Value = 1132
$ws.Range("N126").Value = -8835.0002
$ws.Range("H136").Value = 2029.746
$ws.Range("I136").Value = 1966.317
$ws.Range("J136").Value = 2147.9546
$ws.Range("K136").Value = 5898.951
$ws.Range("L136").Value = 6443.8638
$ws.Range("M136").Value = -3348.951
$ws.Range("N136").Value = -11543.8638
